$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.123.87'
$ws.Range("E2").Value = '  +0.40%  '

$ws.Range("D3").Value = '3.120.66'
$ws.Range("E3").Value = '  +0.65%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.22'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.97%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.522'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.26%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.47'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.74%  '

$ws.Range("E10").Value = '  -0.05%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.480'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.64%  '

$ws.Range("E12").Value = '  +0.03%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.08'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.34%  '

$ws.Range("E14").Value = '  -1.65%  '

$ws.Range("D15").Value = '3.641.47'
$ws.Range("E15").Value = '  +0.75%  '

$ws.Range("D16").Value = '67.124.29'
$ws.Range("E16").Value = '  +0.39%  '

$ws.Range("E17").Value = '  -0.54%  '

$ws.Range("D19").Value = '3.122.45'
$ws.Range("E19").Value = '  +0.64%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '490.11'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.92'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.707'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.99%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '84.20'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.35%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.20'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.30'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.65%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.46'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.99%  '

$ws.Range("E27").Value = '  -0.01%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.92'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.55%  '

$ws.Range("E29").Value = '  -1.75%  '

$ws.Range("E30").Value = '  -0.52%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '28.59'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.89%  '

$ws.Range("E32").Value = '  -0.76%  '

$ws.Range("D33").Value = '0.0₃0950'
$ws.Range("E33").Value = '  -5.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.06%  '

$ws.Range("E35").Value = '  -0.04%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.975'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.84%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '47.26'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.68%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.05'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.83%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.310'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.64%  '

$ws.Range("E40").Value = '  +1.69%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.52'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.56%  '

$ws.Range("D42").Value = '2.822.94'
$ws.Range("E42").Value = '  -0.56%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '385.35'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.14%  '

$ws.Range("E44").Value = '  -7.03%  '

$ws.Range("E45").Value = '  -2.14%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '135.62'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.50%  '

$ws.Range("E48").Value = '  -0.05%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.20'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.32%  '

$ws.Range("E50").Value = '  -0.61%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.75'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.78%  '
